# Updated cryptos list (price & 1h-volume refresh), matching a GitHub
# Actions scheduled run. Rows 38 and 39 swapped their coin identity
# (TheGraph <-> Dai) in addition to value refreshes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 38 & 39: coin/link swap plus new price & volume figures ---
$ws.Range("B38").Value = "Dai"
$ws.Range("C38").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.00"
$ws.Range("E38").Value = "  +0.15%  "

$ws.Range("B39").Value = "TheGraph"
$ws.Range("C39").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.398"
$ws.Range("E39").Value = "  -5.66%  "

# --- Remaining rows: Price (D) / Volume(1h) (E) refresh ---
$ws.Range("D2").Value = "66.519.04"
$ws.Range("E2").Value = "  -4.00%  "
$ws.Range("D3").Value = "3.569.09"
$ws.Range("E3").Value = "  -4.58%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "583.45"
$ws.Range("E5").Value = "  -5.02%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "186.65"
$ws.Range("E6").Value = "  -0.62%  "
$ws.Range("D7").Value = "3.563.87"
$ws.Range("E7").Value = "  -4.56%  "
$ws.Range("E8").Value = "  -4.40%  "
$ws.Range("E9").Value = "  +0.16%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.670"
$ws.Range("E10").Value = "  -7.43%  "
$ws.Range("E11").Value = "  -10.21%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "53.30"
$ws.Range("E12").Value = "  -7.13%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000262"
$ws.Range("E13").Value = "  -11.14%  "
$ws.Range("E14").Value = "  -8.20%  "
$ws.Range("D15").Value = "4.132.66"
$ws.Range("E15").Value = "  -4.38%  "
$ws.Range("D16").Value = "3.566.50"
$ws.Range("E16").Value = "  -4.47%  "
$ws.Range("E17").Value = "  -0.86%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "18.36"
$ws.Range("E18").Value = "  -5.58%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.23"
$ws.Range("E19").Value = "  -6.89%  "
$ws.Range("D20").Value = "66.346.88"
$ws.Range("E20").Value = "  -3.89%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.06"
$ws.Range("E21").Value = "  -7.72%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "395.68"
$ws.Range("E22").Value = "  -4.63%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.35"
$ws.Range("E23").Value = "  -6.48%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "86.07"
$ws.Range("E24").Value = "  -4.04%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.30"
$ws.Range("E25").Value = "  +1.61%  "
$ws.Range("E26").Value = "  -5.57%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.50"
$ws.Range("E27").Value = "  -3.88%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.04"
$ws.Range("E28").Value = "  -0.21%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.54"
$ws.Range("E29").Value = "  -7.10%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.96"
$ws.Range("E30").Value = "  -7.98%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "31.16"
$ws.Range("E31").Value = "  -6.89%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.09"
$ws.Range("E32").Value = "  -4.29%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "12.22"
$ws.Range("E33").Value = "  -4.78%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "622.63"
$ws.Range("E34").Value = "  -0.90%  "
$ws.Range("E35").Value = "  -8.66%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "63.45"
$ws.Range("E36").Value = "  -4.42%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "41.51"
$ws.Range("E37").Value = "  -8.04%  "
$ws.Range("E40").Value = "  -9.69%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.132"
$ws.Range("E41").Value = "  -6.72%  "
$ws.Range("E42").Value = "  -0.04%  "
$ws.Range("D43").Value = "3.010.71"
$ws.Range("E43").Value = "  +5.41%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.83"
$ws.Range("E44").Value = "  -8.25%  "
$ws.Range("E45").Value = "  -5.07%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0409"
$ws.Range("E46").Value = "  -8.56%  "
$ws.Range("E47").Value = "  -0.06%  "
$ws.Range("E48").Value = "  -7.25%  "
$ws.Range("E49").Value = "  -6.86%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "137.28"
$ws.Range("E50").Value = "  -3.82%  "
$ws.Range("E51").Value = "  -1.46%  "
